# Insert two new data rows before row 460 (shifting existing rows 460-489 down to 462-491),
# then populate the two newly-inserted rows with the new "Camote"/"Paine" records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 460.
$ws.Range("A460:R461").EntireRow.Insert()

# New row 460: Camote, 1a nueva(o)
$ws.Range("A460").Value = 10
$ws.Range("B460").Value = "Vega Modelo de Temuco"
$ws.Range("C460").Value = "La Araucanía"
$ws.Range("D460").Value = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D460").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E460").Value = 9
$ws.Range("F460").Value = 100112045
$ws.Range("G460").Value = "Zapallo"
$ws.Range("H460").Value = "Camote"
$ws.Range("I460").Value = "1a nueva(o)"
$ws.Range("J460").Value = 300
$ws.Range("K460").Value = 450
$ws.Range("L460").Value = 450
$ws.Range("M460").Value = 450
$ws.Range("N460").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O460").Value = "Región del Maule"
$ws.Range("P460").Value = 450
$ws.Range("Q460").Value = 1
$ws.Range("R460").Value = "Hortaliza"

# New row 461: Paine, 1a nueva(o)
$ws.Range("A461").Value = 10
$ws.Range("B461").Value = "Vega Modelo de Temuco"
$ws.Range("C461").Value = "La Araucanía"
$ws.Range("D461").Value = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D461").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E461").Value = 9
$ws.Range("F461").Value = 100112045
$ws.Range("G461").Value = "Zapallo"
$ws.Range("H461").Value = "Paine"
$ws.Range("I461").Value = "1a nueva(o)"
$ws.Range("J461").Value = 1000
$ws.Range("K461").Value = 250
$ws.Range("L461").Value = 250
$ws.Range("M461").Value = 250
$ws.Range("N461").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O461").Value = "Región del Maule"
$ws.Range("P461").Value = 250
$ws.Range("Q461").Value = 1
$ws.Range("R461").Value = "Hortaliza"
